# Universidad.xlsx: merge in the new "Capacidad" (capacity) column on the
# "Aulas" sheet and make "Aulas" the active/selected sheet, per the commit
# "Excels unificados ... Se combinaron los archivos anteriores de Edificios
# y Aulas, en el excel 'Universidad.xlsx'".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aulas")

# Insert a new column before the old "Edificio" column (B) and fill its
# header with the new field name. Excel shifts every cell, column width,
# and the data-validation range (B2:B1004 -> C2:C1004) automatically.
$ws.Columns.Item(2).Insert() | Out-Null
$ws.Range("B1").Value = "Capacidad"

# Normalize the (already-broken) #NAME? literal stored in the "Edificio"
# column's validation formula.
$ws.Range("C2:C1004").Validation.Formula1 = "=#¿nombre?"

# "Aulas" becomes the active sheet/tab (previously "Materias" was active),
# with the selection resting on the first empty "Capacidad" cell (B4).
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
